$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Existing column B (dbExcel / Neo4j file
# name) and column C (WebExcel / web file name) shift right to C and D.
$ws.Columns("B:B").Insert()

# New column B header + value ("StatQuery" query used for the stat bar).
$ws.Range("B1").Value = "StatQuery"

$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Chinese Shar-Pei']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").Value = $statQuery

# Match column B's width/wrap formatting to column A.
$ws.Columns("B:B").ColumnWidth = 75
$ws.Range("B2").WrapText = $true

# Selection now covers the whole new StatQuery column.
$ws.Columns("B:B").Select()

$wb.Save()
